$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $basicAmount = $ws.Cells.Item($row, 22).Value2  # Column V = BASIC_AMOUNT
    $ws.Cells.Item($row, 31).Value = $basicAmount  # Column AE = ADJUSTMENT_AMOUNT

    # Columns AK (37) through AO (41) -> set text "0"
    for ($col = 37; $col -le 41; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = "0"
    }
}
